$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text/coin-name updates (pure text cells, safe to set directly)
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E43").Value = "42KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"

# Apply price updates (Column D) as text, preserving inline-string/text semantics
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($c in $priceCells) { $ws.Range($c).NumberFormat = "@" }
$ws.Range("D2").Value = "243.97"
$ws.Range("D3").Value = "23.99"
$ws.Range("D4").Value = "5.265"
$ws.Range("D5").Value = "0.05821"
$ws.Range("D6").Value = "6.460"
$ws.Range("D7").Value = "3.230"
$ws.Range("D8").Value = "0.8087"
$ws.Range("D9").Value = "0.8855"
$ws.Range("D11").Value = "0.07140"
$ws.Range("D13").Value = "0.03041"
$ws.Range("D14").Value = "0.09339"
$ws.Range("D15").Value = "3.836"
$ws.Range("D16").Value = "0.001541"
$ws.Range("D17").Value = "0.04723"
$ws.Range("D18").Value = "0.0006022"
$ws.Range("D19").Value = "0.006201"
$ws.Range("D20").Value = "0.001261"
$ws.Range("D21").Value = "0.004078"
$ws.Range("D22").Value = "0.00008708"
$ws.Range("D23").Value = "3.542"
$ws.Range("D24").Value = "2.159"
$ws.Range("D25").Value = "0.3184"
$ws.Range("D26").Value = "0.1314"
$ws.Range("D41").Value = "0.1055"
$ws.Range("D42").Value = "0.002491"
$ws.Range("D43").Value = "0.006266"
$ws.Range("D44").Value = "0.007226"
$ws.Range("D45").Value = "0.00005335"
$ws.Range("D47").Value = "0.5352"
$ws.Range("D48").Value = "0.005006"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.0002001"
foreach ($c in $priceCells) { $ws.Range($c).Style = "Normal" }
